# #59: fixed data provider RAM values for validation
#
# The GSMArena sheet's "ram" column for the Galaxy S10+ row incorrectly
# read "12GB RAM" - that model actually shipped as either an 8GB or
# 12GB variant, matching the "8/12GB RAM" value already used for the
# Galaxy S10 row. Correct the data-provider fixture accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GSMArena")

# Row 2 = Galaxy S10+, column G = ram
$ws.Range("G2").Value = "8/12GB RAM"
